$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert a new column at L (12) for "Section Code" ------------------
$ws.Columns.Item(12).Insert()

# --- 2. Sample data row (row 2) + header (row 1) ---------------------------
# Order matters: this is the order new shared strings get appended in.
$ws.Range("N2").Value2 = "2023-2024"
$ws.Range("K2").Value2 = "BSCPE"
$ws.Range("M2").Value2 = "Bachelor's of Science in Computer Engineering"
$ws.Range("F2").Value2 = "sample_aemail@gmail.com"
$ws.Range("L1").Value2 = "Section Code"
$ws.Range("L2").Value2 = "3P"

# --- 4. Shrink the old placeholder columns ---------------------------------
# F column: only row 3 keeps a (re-styled) placeholder, rows 4-21 are cleared
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4:F21").Clear()

# Q column (old P, after the column insert) placeholders fully removed
$ws.Range("Q3:Q21").Clear()

# --- 5. New placeholder cells in column H (Mobile Number) ------------------
$ws.Range("H3:H10").Style = "Normal"
$ws.Range("H13:H92").Style = "Normal"

# --- 6. Selection / view ----------------------------------------------------
$ws.Range("J1:J2").Select()
